$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.006.26"
$ws.Range("E2").Value = "  -2.25%  "
$ws.Range("D3").Value = "2.907.04"
$ws.Range("E3").Value = "  -2.33%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'371.37"
$ws.Range("E5").Value = "  +4.78%  "
$ws.Range("D6").Value = "'101.72"
$ws.Range("E6").Value = "  -5.15%  "
$ws.Range("D7").Value = "'0.541"
$ws.Range("E7").Value = "  -3.91%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.586"
$ws.Range("E9").Value = "  -4.69%  "
$ws.Range("D10").Value = "'36.78"
$ws.Range("E10").Value = "  -3.89%  "
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("E12").Value = "  -2.67%  "
$ws.Range("D13").Value = "'18.26"
$ws.Range("E13").Value = "  -5.03%  "
$ws.Range("D14").Value = "3.363.12"
$ws.Range("E14").Value = "  -2.34%  "
$ws.Range("D15").Value = "'7.34"
$ws.Range("E15").Value = "  -3.83%  "
$ws.Range("D16").Value = "2.912.63"
$ws.Range("E16").Value = "  -2.14%  "
$ws.Range("D17").Value = "'0.923"
$ws.Range("E17").Value = "  -7.78%  "
$ws.Range("D18").Value = "50.978.36"
$ws.Range("E18").Value = "  -2.32%  "
$ws.Range("D19").Value = "'3.24"
$ws.Range("E19").Value = "  -6.89%  "
$ws.Range("E20").Value = "  -4.00%  "
$ws.Range("D21").Value = "'12.90"
$ws.Range("E21").Value = "  -5.13%  "
$ws.Range("E22").Value = "  -3.55%  "
$ws.Range("D23").Value = "'68.05"
$ws.Range("E23").Value = "  -2.20%  "
$ws.Range("D24").Value = "'258.80"
$ws.Range("E24").Value = "  -1.95%  "
$ws.Range("E25").Value = "  -2.46%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "'0.168"
$ws.Range("E26").Value = "  -6.07%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").Value = "'4.10"
$ws.Range("E28").Value = "  -4.46%  "
$ws.Range("E29").Value = "  -4.61%  "
$ws.Range("E30").Value = "  -7.31%  "
$ws.Range("E31").Value = "  -7.03%  "
$ws.Range("D32").Value = "'6.24"
$ws.Range("E32").Value = "  +1.71%  "
$ws.Range("E33").Value = "  -4.33%  "
$ws.Range("E34").Value = "  -2.82%  "
$ws.Range("D35").Value = "'51.27"
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("E36").Value = "  -6.24%  "
$ws.Range("E38").Value = "  -4.97%  "
$ws.Range("E39").Value = "  -6.96%  "
$ws.Range("D40").Value = "'16.96"
$ws.Range("E40").Value = "  -5.34%  "
$ws.Range("E41").Value = "  -5.36%  "
$ws.Range("E42").Value = "  -6.79%  "
$ws.Range("E43").Value = "  -3.91%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "'119.50"
$ws.Range("E44").Value = "  -1.68%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'21.83"
$ws.Range("E45").Value = "  -4.05%  "
$ws.Range("E46").Value = "  -2.17%  "
$ws.Range("D47").Value = "2.014.51"
$ws.Range("E47").Value = "  -4.89%  "
$ws.Range("E48").Value = "  -1.91%  "
$ws.Range("E49").Value = "  -7.45%  "
$ws.Range("D50").Value = "3.193.93"
$ws.Range("E50").Value = "  -2.27%  "
$ws.Range("E51").Value = "  -1.58%  "
